$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("boson" shifts from E to F, etc.),
# which also shifts the formulas in (old) J/K/L -> (new) K/L/M and adjusts
# their G/H/I/J references automatically.
$ws.Columns("E:E").Insert()

# Fill the new column E with header "pt_max" and value 50 for every data row.
$ws.Range("E1").Value = "pt_max"
$ws.Range("E2:E12").Value = 50

# Match the author's final selection shown in the saved workbook.
$ws.Range("E17").Select()
